# Auto-generated Excel COM-interop script
# Applies updated market-price/profit figures scraped by the scheduled runner
# to the Leve profit tables across multiple job sheets.

$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H29").Value = 614.5714
$ws.Range("I29").Value = 260.4
$ws.Range("K29").Value = 781.1999999999999
$ws.Range("M29").Value = -500.1999999999999
$ws.Range("H31").Value = 999
$ws.Range("I31").Value = 999
$ws.Range("K31").Value = 2997
$ws.Range("M31").Value = -2767
$ws.Range("H38").Value = 2375.1667
$ws.Range("I38").Value = 229.5
$ws.Range("K38").Value = 688.5
$ws.Range("M38").Value = -316.5
$ws.Range("H55").Value = 409.42856
$ws.Range("I55").Value = 394
$ws.Range("J55").Value = 430
$ws.Range("K55").Value = 394
$ws.Range("L55").Value = 430
$ws.Range("M55").Value = -180
$ws.Range("N55").Value = -858
$ws.Range("H137").Value = 2994
$ws.Range("I137").Value = 2994
$ws.Range("K137").Value = 8982
$ws.Range("M137").Value = -6432

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 227.42857
$ws.Range("I4").Value = 215.33333
$ws.Range("J4").Value = 300
$ws.Range("K4").Value = 215.33333
$ws.Range("L4").Value = 300
$ws.Range("M4").Value = -99.33332999999999
$ws.Range("N4").Value = -532
$ws.Range("H5").Value = 200.42857
$ws.Range("I5").Value = 191
$ws.Range("J5").Value = 213
$ws.Range("K5").Value = 191
$ws.Range("L5").Value = 213
$ws.Range("M5").Value = -79
$ws.Range("N5").Value = -437
$ws.Range("H63").Value = 1850
$ws.Range("I63").Value = 1700
$ws.Range("J63").Value = 2000
$ws.Range("K63").Value = 1700
$ws.Range("L63").Value = 2000
$ws.Range("M63").Value = -1014
$ws.Range("N63").Value = -3372
$ws.Range("H66").Value = 1850
$ws.Range("I66").Value = 1700
$ws.Range("J66").Value = 2000
$ws.Range("K66").Value = 8500
$ws.Range("L66").Value = 10000
$ws.Range("M66").Value = -5068
$ws.Range("N66").Value = -16864
$ws.Range("H101").Value = 33333.668
$ws.Range("J101").Value = 33333.668
$ws.Range("L101").Value = 33333.668
$ws.Range("N101").Value = -39823.668

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 200.42857
$ws.Range("I4").Value = 191
$ws.Range("J4").Value = 213
$ws.Range("K4").Value = 191
$ws.Range("L4").Value = 213
$ws.Range("M4").Value = -76
$ws.Range("N4").Value = -443
$ws.Range("H20").Value = 5868
$ws.Range("I20").Value = 5827
$ws.Range("J20").Value = 5950
$ws.Range("K20").Value = 5827
$ws.Range("L20").Value = 5950
$ws.Range("M20").Value = -5580
$ws.Range("N20").Value = -6444
$ws.Range("H44").Value = 2999.5
$ws.Range("J44").Value = 2999.5
$ws.Range("L44").Value = 2999.5
$ws.Range("N44").Value = -3993.5
$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()  # remove stale N62 value (-31372)
$ws.Range("H63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").ClearContents()  # remove stale N63 value (-51643)
$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()  # remove stale N65 value (-96864)
$ws.Range("H66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").ClearContents()  # remove stale N66 value (-157677)
$ws.Range("H134").Value = 1375.7059
$ws.Range("I134").Value = 1360.8125
$ws.Range("J134").Value = 1614
$ws.Range("K134").Value = 4082.4375
$ws.Range("L134").Value = 4842
$ws.Range("M134").Value = -1547.4375
$ws.Range("N134").Value = -9912

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 388.875
$ws.Range("I22").Value = 388.875
$ws.Range("K22").Value = 388.875
$ws.Range("M22").Value = -38.875
$ws.Range("H31").Value = 3341.25
$ws.Range("I31").Value = 1965
$ws.Range("J31").Value = 4167
$ws.Range("K31").Value = 1965
$ws.Range("L31").Value = 4167
$ws.Range("M31").Value = -1670
$ws.Range("N31").Value = -4757
$ws.Range("H34").Value = 3341.25
$ws.Range("I34").Value = 1965
$ws.Range("J34").Value = 4167
$ws.Range("K34").Value = 1965
$ws.Range("L34").Value = 4167
$ws.Range("M34").Value = -1763
$ws.Range("N34").Value = -4571
$ws.Range("H132").Value = 2765.476
$ws.Range("I132").Value = 2151.4707
$ws.Range("K132").Value = 6454.4121
$ws.Range("M132").Value = -3924.4121
$ws.Range("H134").Value = 1995.25
$ws.Range("I134").Value = 1995.3334
$ws.Range("J134").Value = 1995
$ws.Range("K134").Value = 5986.0002
$ws.Range("L134").Value = 5985
$ws.Range("M134").Value = -3451.0002
$ws.Range("N134").Value = -11055

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 395.66666
$ws.Range("I12").Value = 342.5
$ws.Range("K12").Value = 1027.5
$ws.Range("M12").Value = -854.5
$ws.Range("H121").Value = 1598.5264
$ws.Range("J121").Value = 1663.0588
$ws.Range("L121").Value = 4989.1764
$ws.Range("N121").Value = -7609.1764
$ws.Range("H131").Value = 3000
$ws.Range("J131").Value = 3000
$ws.Range("L131").Value = 9000
$ws.Range("N131").Value = -19080

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 905.5
$ws.Range("I113").Value = 511
$ws.Range("J113").Value = 1300
$ws.Range("K113").Value = 511
$ws.Range("L113").Value = 1300
$ws.Range("M113").Value = 1659
$ws.Range("N113").Value = -5640

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 5362.5
$ws.Range("I16").Value = 3943.75
$ws.Range("K16").Value = 3943.75
$ws.Range("M16").Value = -3773.75
$ws.Range("H40").Value = 3741.3333
$ws.Range("I40").Value = 3741.3333
$ws.Range("K40").Value = 3741.3333
$ws.Range("M40").Value = -3605.3333
$ws.Range("H46").Value = 1395
$ws.Range("I46").Value = 1395
$ws.Range("K46").Value = 1395
$ws.Range("M46").Value = -1207
$ws.Range("H58").Value = 50000
$ws.Range("J58").Value = 50000
$ws.Range("L58").Value = 50000
$ws.Range("N58").Value = -50520
$ws.Range("H75").Value = 25000
$ws.Range("I75").Value = 25000
$ws.Range("K75").Value = 25000
$ws.Range("M75").Value = -24064
$ws.Range("H78").Value = 25000
$ws.Range("I78").Value = 25000
$ws.Range("K78").Value = 75000
$ws.Range("M78").Value = -70320
$ws.Range("H136").Value = 84581
$ws.Range("I136").Value = 5507.4
$ws.Range("K136").Value = 16522.2
$ws.Range("M136").Value = -13972.2
$ws.Range("H140").Value = 68000
$ws.Range("J140").Value = 68000
$ws.Range("L140").Value = 68000
$ws.Range("N140").Value = -78360
$ws.Range("H141").Value = 99999
$ws.Range("J141").Value = 99999
$ws.Range("L141").Value = 99999
$ws.Range("N141").Value = -110359

